# Update "想去人数" (want-to-go count) figures to the latest scrape results.
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$wsExpo = $wb.Worksheets.Item(1)
$wsExpo.Range("F2").Value = 5392
$wsExpo.Range("F4").Value = 635
$wsExpo.Range("F5").Value = 306
$wsExpo.Range("F6").Value = 813
$wsExpo.Range("F7").Value = 8
$wsExpo.Range("F8").Value = 329

# Sheet 2: 演出 (Shows)
$wsShow = $wb.Worksheets.Item(2)
$wsShow.Range("F2").Value = 41
$wsShow.Range("F3").Value = 14

# Sheet 4: 全部类型 (All types / combined sheet)
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F2").Value = 5392
$wsAll.Range("F4").Value = 635
$wsAll.Range("F5").Value = 306
$wsAll.Range("F6").Value = 813
$wsAll.Range("F7").Value = 8
$wsAll.Range("F8").Value = 41
$wsAll.Range("F9").Value = 329
$wsAll.Range("F11").Value = 14
